# Re-position the four lower graph rows on the "Graphs" sheet so the picture
# groups sit closer together (old anchor rows 30/60/90 -> new anchor rows
# 26/52/78; the top group anchored at row 0 is left untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row height is uniform (defaultRowHeight) across this sheet, so the pixel/point
# offset for a given 0-based row index is simply row * rowHeight.
$rowHeight = $ws.Rows.Item(1).RowHeight

# Map of old anchor row -> new anchor row, taken from the target layout.
$rowMap = @{ 30 = 26; 60 = 52; 90 = 78 }

for ($i = 1; $i -le $ws.Shapes.Count; $i++) {
    $shp = $ws.Shapes.Item($i)

    # Which 0-based anchor row is this picture currently sitting on?
    $currentRow = [Math]::Round($shp.Top / $rowHeight)

    if ($rowMap.ContainsKey($currentRow)) {
        $newRow = $rowMap[$currentRow]
        $shp.Top = $newRow * $rowHeight
    }
}
